$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b2 = @'
buenas tardes al generar la información de la cuenta 6960 de as400 está trayendo la información en el documento 3 por error un punto y un cero, adicionando esto al nit o cedula lo cual no permite que el aplicativo de conciliaciones bancarias conban realice el proceso de cruces automáticos. fecha movimiento número cuenta cuenta puc documento3 descripción valor partida tipo partida origen 20/01/2025 4521016960 711140511301 900336004.0 9003360047 administradora colombiana de1111 (62,957,951.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 800113672.0 8001136727 gobierno departamental del t1111 (45,878,783.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 72095111.0 72095111 rodriguez rodriguez pedro anto3167011610 (5,671,899.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 11385211.0 11385211 rodriguez ortiz jose luis 53429415035 (3,578,900.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 890701933.0 8907019334 municipio de melgar ( fonpet1111 (3,446,786.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 52014612.0 52014612 neira ortiz diana rocio 80678439669 (3,206,876.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 51949990.0 51949990 quiroga robayo claudia jasblei35500016351 (2,336,020.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 30309786.0 30309786 ardila pena sandra lucia 91264317652 (2,315,174.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 32705158.0 32705158 silva mejia orietta maria 027900057137 (796,167.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 32698783.0 32698783 silva de la rosa rosmery de je5109029019 (580,644.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 43512185.0 43512185 correa torres mabel insnury 21818699112 (361,925.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 63326901.0 63326901 espitia jimenez nohora smith 3204777598 (243,679.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 71632856.0 71632856 villa bustamante juan carlos 00943190895 (163,711.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 19251557.0 19251557 barbosa estepa reinaldo 24095498803 (29,529.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 70577749.0 70577749 herrera rua javier de jesus 10340868006 (24,319.00) nota débito no extracto libros 20/01/2025 4521016960 711140511301 51783666.0 51783666 piarpuzan ravelo maria de los 20770028944 (11,028.00) nota débito no extracto libros cordialmente, nota: esta comunicación de partidas conciliatorias no constituye soporte para la elaboración y aplicación de ajustes. el área dueña del proceso será la responsable del análisis, gestión y documentación de forma independiente a esta información, para la depuración de las partidas conciliatorias. la información transmitida está destinada al uso exclusivo de las personas a quienes esté dirigida y puede contener información confidencial. toda revisión, retransmisión, diseminación u otro uso de esta información por personas o entidades distintas del destinatario está prohibida. si usted ha recibido esta comunicación por error, por favor haga caso omiso de su contenido y notifíquenos de inmediato a jurídico y cumplimiento colombia. the information transmitted is intended only for the use of the individual to whom it is addressed and may contain information that is confidential or privileged. any review, retransmission, dissemination or other use of, or taking of any action in reliance upon, this information by individuals or entities other than the intended recipient is prohibited. if you have received this communication by mistake, please disregard its contents and notify us immediately to legal and compliance colombia.
'@
$ws.Range("B2").Value = $b2

$c2 = @'
[-0.7204535  -0.65437222 -0.02899576 ...  0.15441492 -0.26677454
  0.6289354 ]
'@
$ws.Range("C2").Value = $c2

$b3 = @'
buenas tardes, por favor su ayuda con la restauración de journal en ambiente prd en libdatple de: cuear1 cuear9 para el rango de fechas de 2024/10/15 a 2024/10/30 cordialmente, andres felipe guerra g. ingeniero backend calle 67 nº. 7 – 94 piso 18 aguerra@colfondos.com.co bogotá d.c. colombia visítanos en www.colfondos.com.co síguenos en: la información transmitida está destinada al uso exclusivo de las personas a quienes esté dirigida y puede contener información confidencial. toda revisión, retransmisión, diseminación u otro uso de esta información por personas o entidades distintas del destinatario está prohibida. si usted ha recibido esta comunicación por error, por favor haga caso omiso de su contenido y notifíquenos de inmediato a jurídico y cumplimiento colombia. the information transmitted is intended only for the use of the individual to whom it is addressed and may contain information that is confidential or privileged. any review, retransmission, dissemination or other use of, or taking of any action in reliance upon, this information by individuals or entities other than the intended recipient is prohibited. if you have received this communication by mistake, please disregard its contents and notify us immediately to legal and compliance colombia.
'@
$ws.Range("B3").Value = $b3

$c3 = @'
[-0.54384172  0.00649924 -0.0143388  ... -0.24499919 -0.53908944
  0.1277236 ]
'@
$ws.Range("C3").Value = $c3

$b4 = @'
datos nombre usuario sede: nombre de la aplicación/producto y no. de versión: pension voluntaria descripción de la solicitud/ error/ inconveniente: se solicitan journals de la maquina de produccion de la libreria libdatple de - cuearc para el rango de fechas 2024/12/25 hasta 2024/12/30 - movarc para la fecha 2024/12/26​
'@
$ws.Range("B4").Value = $b4

$c4 = @'
[-0.30551696 -0.34293187 -0.41503149 ... -0.34063587 -0.08820747
 -0.03934072]
'@
$ws.Range("C4").Value = $c4

$b5 = @'
buen día solicito de su amable colaboración con el error que me genera al momento de dar f10 para que queden los datos actualizados quedo atenta a sus comentarios laura nikole ibañez garcia representante de servicio al cliente of. puente aranda libanez@colfondos.com.co bogotá d.c. colombia ¡accede al código qr, conoce nuestras ofertas laborales y compártelas! la información transmitida está destinada al uso exclusivo de las personas a quienes esté dirigida y puede contener información confidencial. toda revisión, retransmisión, diseminación u otro uso de esta información por personas o entidades distintas del destinatario está prohibida. si usted ha recibido esta comunicación por error, por favor haga caso omiso de su contenido y notifíquenos de inmediato a jurídico y cumplimiento colombia. the information transmitted is intended only for the use of the individual to whom it is addressed and may contain information that is confidential or privileged. any review, retransmission, dissemination or other use of, or taking of any action in reliance upon, this information by individuals or entities other than the intended recipient is prohibited. if you have received this communication by mistake, please disregard its contents and notify us immediately to legal and compliance colombia.
'@
$ws.Range("B5").Value = $b5

$c5 = @'
[-0.52957761  0.11498886 -0.0226076  ...  0.04986332 -0.27396405
  0.09733911]
'@
$ws.Range("C5").Value = $c5

$b6 = @'
buen día, por favor ejecutar aplicación m23 en ambiente uat. nota: @buzon mesa de servicio por favor crear caso y asignar a @jose eden aroca ibarra atentamente, david steven martinez vargas ingeniero de aplicaciones calle 67 nº. 7 – 94 teléfono: 3765155 ext 10528 dmartinez@colfondos.com.co bogotá d.c. colombia la información transmitida está destinada al uso exclusivo de las personas a quienes esté dirigida y puede contener información confidencial. toda revisión, retransmisión, diseminación u otro uso de esta información por personas o entidades distintas del destinatario está prohibida. si usted ha recibido esta comunicación por error, por favor haga caso omiso de su contenido y notifíquenos de inmediato a jurídico y cumplimiento colombia. the information transmitted is intended only for the use of the individual to whom it is addressed and may contain information that is confidential or privileged. any review, retransmission, dissemination or other use of, or taking of any action in reliance upon, this information by individuals or entities other than the intended recipient is prohibited. if you have received this communication by mistake, please disregard its contents and notify us immediately to legal and compliance colombia. la información transmitida está destinada al uso exclusivo de las personas a quienes esté dirigida y puede contener información confidencial. toda revisión, retransmisión, diseminación u otro uso de esta información por personas o entidades distintas del destinatario está prohibida. si usted ha recibido esta comunicación por error, por favor haga caso omiso de su contenido y notifíquenos de inmediato a jurídico y cumplimiento colombia. the information transmitted is intended only for the use of the individual to whom it is addressed and may contain information that is confidential or privileged. any review, retransmission, dissemination or other use of, or taking of any action in reliance upon, this information by individuals or entities other than the intended recipient is prohibited. if you have received this communication by mistake, please disregard its contents and notify us immediately to legal and compliance colombia.
'@
$ws.Range("B6").Value = $b6

$c6 = @'
[-0.63663912 -0.47595757  0.56054562 ...  0.16122428 -0.40702441
  0.12867706]
'@
$ws.Range("C6").Value = $c6

$b7 = @'
se tiene problemas en ingreso al 400 en ambiente dev - uat, que genera impacto en el proyecto colf0454 para ejecución de pruebas tabla 2 sfc
'@
$ws.Range("B7").Value = $b7

$c7 = @'
[-0.22393948  0.2021974  -0.71701121 ... -0.16164231 -0.14379385
 -0.19330895]
'@
$ws.Range("C7").Value = $c7

$b9 = @'
datos nombre usuario leidy tatiana sánchez castrosede: trabajo en casanombre de la aplicación/producto y no. de versión: fallas en as400 descripción de la solicitud/ error/ inconveniente: usuaria se comunica para solicitar que se revise el equipo dado que no le permite abrir as400 y bizagi.
'@
$ws.Range("B9").Value = $b9

$c9 = @'
[-0.37004113 -0.67240089 -0.42205867 ... -0.11446803  0.04496152
 -0.7594046 ]
'@
$ws.Range("C9").Value = $c9

$b10 = @'
datos nombre usuario eduardo alba salazar sede: n/a nombre de la aplicación/producto y no. de versión: correo tecnoquímicas nit 890300466descripción de la solicitud/ error/ inconveniente: buenas tardes me ayudan radicando la siguiente llamada soporte aplicaciones - as400 core - pv- po -ce, cliente único espacio para diligenciar por usuario área operativa que reporta (acreditación - pagos - recaudo, etc.) afiliaciones digital proceso (cobranzas, pqr´s, historia laboral, afiliaciones, etc.) afiliaciones coordinador asignado eduardo alba salazar ruta en as400 y usuario de ejecución usuario: ea10829 soporte de ejecución - pantallazos o evidencias se requiere revisar el servicio de consulta de datos empleador, cuando se digita el nit no esta cargando el email que tiene la empresa diligenciada en la base de datos observaciones del usuario ¿qué resultado esperaba? al consultar el nit debe traer toda la información de la empresa según el servicio que existe y que venia funcionando cordialmente, nota: utiliza el siguiente buzón para la centralización de consultas, fallas y temas que tienen que ver con el proceso de afiliación digital para pensiones obligatorias, cesantías, pensiones voluntarias y autorizador de trámites, buzonsoporteafiliacion@colfondos.com.co
'@
$ws.Range("B10").Value = $b10

$c10 = @'
[-0.06284009 -0.35354272 -0.01566264 ... -0.27085856 -0.66875726
 -0.06565149]
'@
$ws.Range("C10").Value = $c10

$b11 = @'
datosnombre usuarioedgar javier duarte forerosede:torre colfondos piso 18nombre de la aplicación/producto y no. de versión:ibm idescripción de la solicitud/ error/ inconveniente:se requiere el cargue del certificado digital en el servidor ibm i para poder realizar el consumo de los siguientes web services de la aplicación web de vigía:-https://vascobog009:4443/api/utility/testapiverificatercero.php-https://vascobog009:4443/api/api.php/vigiamc/restrictivelist/query?token=xkgs+sbmulfwwg8y+a6qhwzvr05kc1jsovjznwlqbgjtavfqdw56detsrfhnektmunrsl01sahnlohrxrwtoww1teu9fszhbvlhpvvr5djk=&source=testapivt&search_string=10100159792&indicator=1&percent=100lo anterior con el fin de poder realizar la consulta de los terceros de beneficiarios de pagos realizados a través de los productos core en las listas restrictivas de vigía, para elanálisis del área de sarlaft para confirmar o rechazar el pago.la información del servidor de vigía donde se encuentran los servicios a consumir es la siguiente:nombre: vascobog009ip: 172.18.169.22puerto: 4443quedo atento al cargue del certificado. muchas gracias.
'@
$ws.Range("B11").Value = $b11

$c11 = @'
[-0.27856547 -0.38790217  0.49086815 ... -0.25003666 -0.59334457
  0.23946011]
'@
$ws.Range("C11").Value = $c11

$b12 = @'
se requiere el usuario qsecofr para depurar los certificados tls versiones obsoletas 1.1 y 1.2 en el almacen de certificados de la particion de desarrollo.
'@
$ws.Range("B12").Value = $b12

$c12 = @'
[-0.3967391  -0.28182188 -0.18796274 ...  0.21426664 -0.25727648
 -0.12916225]
'@
$ws.Range("C12").Value = $c12

$b13 = @'
buen día: señores mesa de servicio solicitamos por favor escalar con quien corresponda, con el fin de que nos colaboren validando lo ocurrido con los siguientes depósitos de bono pensional: se remito caso al área de acreditaciónes, con el fin de que se procediera con la acreditación del bono pensional por valor de $143.208.000 para el afiliado cc 12719495 munoz penaloza ismael de jesus, el deposito fue recibido 30/06/2022 así: se cargo solicitud de acreditación en winscp, bajo rb521024 (adjunto); sin embargo, y aun cuando nos indican que el deposito se encuentra acreditado, al validar por la siguiente ruta el saldo no se encuentra en estado acr: 8.- futura 24 bono pensional 5 control de bono 1 trabajar con bonos por afiliado cc cc 79543263 - duque cruz fabio herney: para este caso fue escalado en alguna oportunidad con ti con el fin de que se validara porque no se permitia la acreditación del mismo; sin embargo, nos indicaron que el mismo tenia solicitud de traspaso marcado en el solar6, la solicitud se elevo con pensiones pero nos indican que el afiliado no cuenta con solicitud activa de traspaso, por cuanto nos urge acreditar el bono. se envio bajo rb420924 (adjunto) y se nos indica que el mismo se encuentra acreditado; sin embargo al consultar la siguiente ruta el estado no es acr: 8.- futura 24 bono pensional 5 control de bono 1 trabajar con bonos por afiliado cc gracias quedo atenta a comentarios zaida yulied pardo chacón analista de bonos pensiones calle 67 nº. 7 – 94 zpardo@colfondos.com.co bogotá d.c. colombia la información transmitida está destinada al uso exclusivo de las personas a quienes esté dirigida y puede contener información confidencial. toda revisión, retransmisión, diseminación u otro uso de esta información por personas o entidades distintas del destinatario está prohibida. si usted ha recibido esta comunicación por error, por favor haga caso omiso de su contenido y notifíquenos de inmediato a jurídico y cumplimiento colombia. the information transmitted is intended only for the use of the individual to whom it is addressed and may contain information that is confidential or privileged. any review, retransmission, dissemination or other use of, or taking of any action in reliance upon, this information by individuals or entities other than the intended recipient is prohibited. if you have received this communication by mistake, please disregard its contents and notify us immediately to legal and compliance colombia. ========================= forwarded mail ========================= message id : sent date : tue, 29 oct 2024 16:40:08 +0000 from : yramos@colfondos.com.co sender mailid : yramos@colfondos.com.co reply to : yramos@colfondos.com.co to : zpardo@colfondos.com.co,dlhernandez@colfondos.com.co cc : mmantilla@colfondos.com.co,dgonzales-empleadoenmision@colfondos.com.co,jgarcia-empleadoenmision@colfondos.com.co subject : re: para acreditar deposito - cc 12719495 munoz penaloza ismael de jesus//rb521024 ------------------------- contents ------------------------- buen día, acreditación de bono gestionado, por favor validar. quedo atenta. cordialmente, de: zaida yulied pardo chacon <zpardo@colfondos.com.co> enviado el: lunes, 28 de octubre de 2024 5:36 p.m. para: yenny carolina ramos triana <yramos@colfondos.com.co>; dlayane eslendy hernandez gonzalez <dlhernandez@colfondos.com.co> cc: maria stella mantilla parada <mmantilla@colfondos.com.co>; dayana rosalba gonzalez amaya <dgonzales-empleadoenmision@colfondos.com.co>; julieth katerin garcia mancera <jgarcia-empleadoenmision@colfondos.com.co> asunto: rv: para acreditar deposito - cc 12719495 munoz penaloza ismael de jesus//rb521024 buen día carolina me ayudas por favor acreditando $ 143.208.000 (30/06/2022) con nacion. nit 1 gracias quedo atenta a comentarios zaida yulied pardo chacón analista de bonos pensiones calle 67 nº. 7 – 94 zpardo@colfondos.com.co bogotá d.c. colombia de: dayana rosalba gonzalez amaya <dgonzales-empleadoenmision@colfondos.com.co> enviado el: lunes, 28 de octubre de 2024 5:25 p.m. para: zaida yulied pardo chacon <zpardo@colfondos.com.co> cc: julieth katerin garcia mancera <jgarcia-empleadoenmision@colfondos.com.co> asunto: para acreditar deposito - cc 12719495 munoz penaloza ismael de jesus//rb521024 cordial saludo, grabación de cupón y archivo de acreditación ok. solicito de tú amable colaboración en aprobar el archivo en asunto para acreditar $ 143.208.000 (30/06/2022) con nacion. nit 1, te informo que validé la base de préstamos de bonos, jurídica, cuentas por cobrar y optima (as400) y existe coincidencia con elcaso, por cuanto no se reporta valor a recuperar. archivo: rb521024 cordialmente dayana gonzalez amaya auxiliar bonos pensionales. dgonzales-empleadoenmision@colfondos.com.co bogotá d.c. colombia de: julieth katerin garcia mancera <jgarcia-e...
'@
$ws.Range("B13").Value = $b13

$c13 = @'
[-0.10371219  0.11457424 -0.29114372 ... -0.28407115 -0.40188116
 -0.06783258]
'@
$ws.Range("C13").Value = $c13

$b14 = @'
buenos días informo que el día de hoy he intentado anular 2 trámites por fallecimiento del afiliado y no está dejando. se intenta dejar en estado anu por fallecimiento, pero sigue quedando en estado sol. mil gracias. 8.- futura 12 beneficios y pensiones db2 - dds 1 nuevo modulo de benficios y pensiones 04 solicitudes de pensión y auxilio funerario 01 trabajar con solicitudes número de cédula 19 anular solicitud de pensión f4 @mue afiliado fallecido f10 para confirmar marcela betancur vasco analista de indemnizaciones calle 67 nº. 7 – 94 dbetancur@colfondos.com.co bogotá d.c. colombia ¡accede al código qr, conoce nuestras ofertas laborales y compártelas! la información transmitida está destinada al uso exclusivo de las personas a quienes esté dirigida y puede contener información confidencial. toda revisión, retransmisión, diseminación u otro uso de esta información por personas o entidades distintas del destinatario está prohibida. si usted ha recibido esta comunicación por error, por favor haga caso omiso de su contenido y notifíquenos de inmediato a jurídico y cumplimiento colombia. the information transmitted is intended only for the use of the individual to whom it is addressed and may contain information that is confidential or privileged. any review, retransmission, dissemination or other use of, or taking of any action in reliance upon, this information by individuals or entities other than the intended recipient is prohibited. if you have received this communication by mistake, please disregard its contents and notify us immediately to legal and compliance colombia.
'@
$ws.Range("B14").Value = $b14

$c14 = @'
[-0.42706642 -0.2462464  -0.25242588 ... -0.15655087 -0.67980111
  0.30507284]
'@
$ws.Range("C14").Value = $c14

$b15 = @'
buenas tardes, solicitamos su amble colaboración con la apertura de un caso con la siguiente información: aplicación: ibm i producción descripción: por favor generar el certificado digital para cargar en el ibm i de producción con el fin de que se pueda realizar el consumo de los web services de vigía para la consulta de listas restrictivas. por favor asignar el caso a @jaime enrique nieto naranjo. muchas gracias por su atención. cordial saludo, edgar javier duarte forero ingeniero de aplicaciones calle 67 nº. 7 - 94 tel: (571) 376 50 66 - (571) 376 51 55 ext. 10417 eduarte@colfondos.com.co bogotá d.c. colombia la información transmitida está destinada al uso exclusivo de las personas a quienes esté dirigida y puede contener información confidencial. toda revisión, retransmisión, diseminación u otro uso de esta información por personas o entidades distintas del destinatario está prohibida. si usted ha recibido esta comunicación por error, por favor haga caso omiso de su contenido y notifíquenos de inmediato a jurídico y cumplimiento colombia. the information transmitted is intended only for the use of the individual to whom it is addressed and may contain information that is confidential or privileged. any review, retransmission, dissemination or other use of, or taking of any action in reliance upon, this information by individuals or entities other than the intended recipient is prohibited. if you have received this communication by mistake, please disregard its contents and notify us immediately to legal and compliance colombia.
'@
$ws.Range("B15").Value = $b15

$c15 = @'
[-0.23526996 -0.6012882  -0.07465767 ...  0.18934099 -0.48558214
  0.35685769]
'@
$ws.Range("C15").Value = $c15

$b16 = @'
datos nombre usuario leda marina bolivar sandovalsede: barranquillanombre de la aplicación/producto y no. de versión: as-400 descripción de la solicitud/ error/ inconveniente: bloqueo de usuario
'@
$ws.Range("B16").Value = $b16

$c16 = @'
[-0.55536401 -0.30159673 -0.84147847 ...  0.2249067   0.04742323
 -0.11696886]
'@
$ws.Range("C16").Value = $c16

